$d = $word.ActiveDocument

# The document has paragraphs whose <w:pPr> carries two <w:pStyle> elements
# (a leftover "Compact" style followed by the real style). Word's object
# model resolves such a paragraph's effective Style to the *last* pStyle
# entry, so only the paragraph that has "Compact" alone (no second entry)
# currently shows up with Style "Compact" through the object model - it is
# the footnote-content paragraph ("... my note ... http://example.com").
#
# Re-assigning that paragraph's Style to "Body Text" triggers Word to
# rewrite/normalize the paragraph properties for the whole document,
# dropping the redundant leading "Compact" pStyle entries everywhere else
# too, while this paragraph's own style becomes "Body Text".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Compact") {
        $p.Style = "Body Text"
    }
}
